$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Jul 11 12:54:16 EDT 2023"
$ws.Range("B3").Value = "Tue Jul 11 12:54:26 EDT 2023"
$ws.Range("B4").Value = "Tue Jul 11 12:54:35 EDT 2023"
$ws.Range("B5").Value = "Tue Jul 11 12:54:45 EDT 2023"
$ws.Range("B6").Value = "Tue Jul 11 12:54:54 EDT 2023"
$ws.Range("B7").Value = "Tue Jul 11 12:55:03 EDT 2023"
